$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# Price cells hold numeric-looking text (e.g. "514.88"); force text storage with a
# temporary NumberFormat so Excel does not auto-convert them to numbers, then restore
# the cell style so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.119.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.046.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.574.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.242.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.048.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.496"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.09%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0893"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0668"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.087.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.649"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.248.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.52%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "259.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0873"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.84%  "
